$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.60814356803894
$ws.Range("B1").Value = 2.472468137741089
$ws.Range("C1").Value = 2.769110441207886
$ws.Range("D1").Value = 3.718555450439453
$ws.Range("E1").Value = 5.173699378967285
